$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

$ws.Range("F4").Value = 75
$ws.Range("G4").Value = 73.2
$ws.Range("H4").Value = 76.7
$ws.Range("F10").Value = 42.2
$ws.Range("G10").Value = 43.3
$ws.Range("H10").Value = 41.3
$ws.Range("F13").Value = 77.955555555555605
$ws.Range("G13").Value = 80.816666666666706
$ws.Range("H13").Value = 75.133333333333397
$ws.Range("G38").Value = 82.290000000000106
$ws.Range("F62").Value = 82.630769230769303
$ws.Range("G62").Value = 85.098717948717905
$ws.Range("H62").Value = 80.394871794871804
$ws.Range("F63").Value = 40.226751592356699
$ws.Range("G63").Value = 39.480254777070101
$ws.Range("H63").Value = 40.6732484076433
$ws.Range("F64").Value = 57.809375000000003
$ws.Range("G64").Value = 55.381250000000001
$ws.Range("H64").Value = 59.578125
$ws.Range("F65").Value = 73.073684210526295
$ws.Range("G65").Value = 72.592105263157904
$ws.Range("H65").Value = 72.842105263157904
$ws.Range("F66").Value = 54.301276595744703
$ws.Range("G66").Value = 54.621702127659603
$ws.Range("H66").Value = 53.857446808510701
$ws.Range("G67").Value = 80.694117647058803
$ws.Range("G68").Value = 86.8333333333333
$ws.Range("F73").Value = 78.892857142857196
$ws.Range("G73").Value = 80.967857142857198
$ws.Range("H73").Value = 76.807142857142907
$ws.Range("F76").Value = 71.599999999999994
$ws.Range("G76").Value = 72.825000000000003
$ws.Range("H76").Value = 70.733333333333405
$ws.Range("F77").Value = 55.193333333333399
$ws.Range("G77").Value = 53.62
$ws.Range("H77").Value = 56.1933333333333
$ws.Range("F78").Value = 8.0346153846153801
$ws.Range("G78").Value = 7.2923076923076904
$ws.Range("F79").Value = 12.5129032258065
$ws.Range("G79").Value = 12.3032258064516
$ws.Range("H79").Value = 12.8403225806452
$ws.Range("F81").Value = 59.1933333333333
$ws.Range("G81").Value = 55.626666666666701
$ws.Range("H81").Value = 60.566666666666698
$ws.Range("F82").Value = 82.348571428571503
$ws.Range("G82").Value = 84.821428571428598
$ws.Range("H82").Value = 80.099999999999994
$ws.Range("F83").Value = 38.2232394366197
$ws.Range("G83").Value = 37.774647887324001
$ws.Range("H83").Value = 38.571830985915497
$ws.Range("G84").Value = 93.657894736842195
$ws.Range("F87").Value = 75.174285714285702
$ws.Range("G87").Value = 75.7628571428572
$ws.Range("H87").Value = 74.948571428571398
$ws.Range("F88").Value = 48.1666666666667
$ws.Range("G88").Value = 47.4166666666667
$ws.Range("H88").Value = 48.466666666666697
$ws.Range("F89").Value = 50.1413043478261
$ws.Range("G89").Value = 47.5717391304348
$ws.Range("H89").Value = 51.05
$ws.Range("F90").Value = 13.34
$ws.Range("G90").Value = 12.7314285714286
$ws.Range("H90").Value = 13.8828571428571
$ws.Range("G91").Value = 92.196296296296396
$ws.Range("H91").Value = 86.044444444444395
$ws.Range("F94").Value = 59.9304347826087
$ws.Range("G94").Value = 57.056521739130403
$ws.Range("H94").Value = 61.9652173913044
$ws.Range("F95").Value = 86.274074074074093
$ws.Range("G95").Value = 89.118518518518499
$ws.Range("H95").Value = 84.085185185185196
$ws.Range("F96").Value = 66.099999999999994
$ws.Range("G96").Value = 67.929411764705904
$ws.Range("H96").Value = 65.276470588235298
$ws.Range("G97").Value = 91.372413793103405
$ws.Range("F98").Value = 80.895652173913007
$ws.Range("G98").Value = 81.547826086956505
$ws.Range("H98").Value = 80.273913043478302
